# Add the "bkk" use case to the config sheet, then make config the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

$ws.Cells.Item(19, 2).Value = "baseYear"
$ws.Cells.Item(19, 3).Value = 2015

$ws.Cells.Item(20, 2).Value = "finalYear"
$ws.Cells.Item(20, 3).Value = 2115

$ws.Cells.Item(21, 2).Value = "zones_filename"
$ws.Cells.Item(21, 3).Value = "map"

$ws.Cells.Item(19, 1).Value = "bkk"
$ws.Cells.Item(20, 1).Value = "bkk"
$ws.Cells.Item(21, 1).Value = "bkk"

$ws.Activate()
$ws.Range("A22").Select()
